$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price / Volume(1h)) per the commit diff.
# Rows 17 and 18 swap Coin/Link (Polygon <-> WrappedEther) along with their price/volume.
# Values are assigned with a leading apostrophe (forces text entry, matching the
# original inline-string cells) and ClearFormats() strips the resulting quote-prefix
# style so the cell keeps its original (default) style index.

$c = $ws.Range('D2')
$c.Value = "'" + '41.299.37'
$c.ClearFormats()
$c = $ws.Range('E2')
$c.Value = "'" + '  -1.54%  '
$c.ClearFormats()
$c = $ws.Range('D3')
$c.Value = "'" + '2.175.85'
$c.ClearFormats()
$c = $ws.Range('E3')
$c.Value = "'" + '  -1.78%  '
$c.ClearFormats()
$c = $ws.Range('E4')
$c.Value = "'" + '  +0.08%  '
$c.ClearFormats()
$c = $ws.Range('D5')
$c.Value = "'" + '237.71'
$c.ClearFormats()
$c = $ws.Range('E5')
$c.Value = "'" + '  -1.31%  '
$c.ClearFormats()
$c = $ws.Range('D6')
$c.Value = "'" + '0.609'
$c.ClearFormats()
$c = $ws.Range('E6')
$c.Value = "'" + '  -2.49%  '
$c.ClearFormats()
$c = $ws.Range('D7')
$c.Value = "'" + '70.00'
$c.ClearFormats()
$c = $ws.Range('E7')
$c.Value = "'" + '  -4.29%  '
$c.ClearFormats()
$c = $ws.Range('D9')
$c.Value = "'" + '0.574'
$c.ClearFormats()
$c = $ws.Range('E9')
$c.Value = "'" + '  -5.16%  '
$c.ClearFormats()
$c = $ws.Range('D10')
$c.Value = "'" + '39.45'
$c.ClearFormats()
$c = $ws.Range('E10')
$c.Value = "'" + '  -7.89%  '
$c.ClearFormats()
$c = $ws.Range('D11')
$c.Value = "'" + '0.0920'
$c.ClearFormats()
$c = $ws.Range('E11')
$c.Value = "'" + '  -3.34%  '
$c.ClearFormats()
$c = $ws.Range('D12')
$c.Value = "'" + '54.49'
$c.ClearFormats()
$c = $ws.Range('E12')
$c.Value = "'" + '  -5.13%  '
$c.ClearFormats()
$c = $ws.Range('E13')
$c.Value = "'" + '  -1.99%  '
$c.ClearFormats()
$c = $ws.Range('D14')
$c.Value = "'" + '6.72'
$c.ClearFormats()
$c = $ws.Range('E14')
$c.Value = "'" + '  -4.89%  '
$c.ClearFormats()
$c = $ws.Range('D15')
$c.Value = "'" + '2.500.28'
$c.ClearFormats()
$c = $ws.Range('E15')
$c.Value = "'" + '  -1.78%  '
$c.ClearFormats()
$c = $ws.Range('D16')
$c.Value = "'" + '14.31'
$c.ClearFormats()
$c = $ws.Range('E16')
$c.Value = "'" + '  +0.67%  '
$c.ClearFormats()
$c = $ws.Range('B17')
$c.Value = "'" + 'WrappedEther'
$c.ClearFormats()
$c = $ws.Range('C17')
$c.Value = "'" + 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c.ClearFormats()
$c = $ws.Range('D17')
$c.Value = "'" + '2.160.19'
$c.ClearFormats()
$c = $ws.Range('E17')
$c.Value = "'" + '  -2.32%  '
$c.ClearFormats()
$c = $ws.Range('B18')
$c.Value = "'" + 'Polygon'
$c.ClearFormats()
$c = $ws.Range('C18')
$c.Value = "'" + 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c.ClearFormats()
$c = $ws.Range('D18')
$c.Value = "'" + '0.793'
$c.ClearFormats()
$c = $ws.Range('E18')
$c.Value = "'" + '  -4.98%  '
$c.ClearFormats()
$c = $ws.Range('D19')
$c.Value = "'" + '41.134.44'
$c.ClearFormats()
$c = $ws.Range('E19')
$c.Value = "'" + '  -1.61%  '
$c.ClearFormats()
$c = $ws.Range('E20')
$c.Value = "'" + '  -6.94%  '
$c.ClearFormats()
$c = $ws.Range('D21')
$c.Value = "'" + '70.57'
$c.ClearFormats()
$c = $ws.Range('E21')
$c.Value = "'" + '  -3.81%  '
$c.ClearFormats()
$c = $ws.Range('D22')
$c.Value = "'" + '5.87'
$c.ClearFormats()
$c = $ws.Range('E22')
$c.Value = "'" + '  -4.63%  '
$c.ClearFormats()
$c = $ws.Range('D23')
$c.Value = "'" + '226.27'
$c.ClearFormats()
$c = $ws.Range('E23')
$c.Value = "'" + '  -1.30%  '
$c.ClearFormats()
$c = $ws.Range('D24')
$c.Value = "'" + '9.38'
$c.ClearFormats()
$c = $ws.Range('E24')
$c.Value = "'" + '  -9.54%  '
$c.ClearFormats()
$c = $ws.Range('E25')
$c.Value = "'" + '  -8.49%  '
$c.ClearFormats()
$c = $ws.Range('E26')
$c.Value = "'" + '  -0.11%  '
$c.ClearFormats()
$c = $ws.Range('D27')
$c.Value = "'" + '10.75'
$c.ClearFormats()
$c = $ws.Range('E27')
$c.Value = "'" + '  -7.47%  '
$c.ClearFormats()
$c = $ws.Range('D28')
$c.Value = "'" + '3.47'
$c.ClearFormats()
$c = $ws.Range('E28')
$c.Value = "'" + '  -3.63%  '
$c.ClearFormats()
$c = $ws.Range('D29')
$c.Value = "'" + '2.20'
$c.ClearFormats()
$c = $ws.Range('E29')
$c.Value = "'" + '  -2.50%  '
$c.ClearFormats()
$c = $ws.Range('E30')
$c.Value = "'" + '  -0.89%  '
$c.ClearFormats()
$c = $ws.Range('D31')
$c.Value = "'" + '167.50'
$c.ClearFormats()
$c = $ws.Range('E31')
$c.Value = "'" + '  +0.58%  '
$c.ClearFormats()
$c = $ws.Range('D32')
$c.Value = "'" + '19.92'
$c.ClearFormats()
$c = $ws.Range('E32')
$c.Value = "'" + '  -3.20%  '
$c.ClearFormats()
$c = $ws.Range('D33')
$c.Value = "'" + '30.45'
$c.ClearFormats()
$c = $ws.Range('E33')
$c.Value = "'" + '  +4.51%  '
$c.ClearFormats()
$c = $ws.Range('D34')
$c.Value = "'" + '0.0762'
$c.ClearFormats()
$c = $ws.Range('E34')
$c.Value = "'" + '  -4.00%  '
$c.ClearFormats()
$c = $ws.Range('D35')
$c.Value = "'" + '5.11'
$c.ClearFormats()
$c = $ws.Range('E35')
$c.Value = "'" + '  -9.37%  '
$c.ClearFormats()
$c = $ws.Range('D36')
$c.Value = "'" + '0.120'
$c.ClearFormats()
$c = $ws.Range('E36')
$c.Value = "'" + '  -3.15%  '
$c.ClearFormats()
$c = $ws.Range('E37')
$c.Value = "'" + '  -7.61%  '
$c.ClearFormats()
$c = $ws.Range('D38')
$c.Value = "'" + '4.09'
$c.ClearFormats()
$c = $ws.Range('E38')
$c.Value = "'" + '  -3.46%  '
$c.ClearFormats()
$c = $ws.Range('D39')
$c.Value = "'" + '0.0282'
$c.ClearFormats()
$c = $ws.Range('E39')
$c.Value = "'" + '  -6.06%  '
$c.ClearFormats()
$c = $ws.Range('E40')
$c.Value = "'" + '  -1.96%  '
$c.ClearFormats()
$c = $ws.Range('D41')
$c.Value = "'" + '11.79'
$c.ClearFormats()
$c = $ws.Range('E41')
$c.Value = "'" + '  -12.10%  '
$c.ClearFormats()
$c = $ws.Range('D42')
$c.Value = "'" + '5.38'
$c.ClearFormats()
$c = $ws.Range('E42')
$c.Value = "'" + '  -3.93%  '
$c.ClearFormats()
$c = $ws.Range('D43')
$c.Value = "'" + '58.89'
$c.ClearFormats()
$c = $ws.Range('E43')
$c.Value = "'" + '  -11.08%  '
$c.ClearFormats()
$c = $ws.Range('D44')
$c.Value = "'" + '0.190'
$c.ClearFormats()
$c = $ws.Range('E44')
$c.Value = "'" + '  -3.61%  '
$c.ClearFormats()
$c = $ws.Range('D45')
$c.Value = "'" + '8.28'
$c.ClearFormats()
$c = $ws.Range('E45')
$c.Value = "'" + '  -4.64%  '
$c.ClearFormats()
$c = $ws.Range('D46')
$c.Value = "'" + '0.0967'
$c.ClearFormats()
$c = $ws.Range('E46')
$c.Value = "'" + '  -3.45%  '
$c.ClearFormats()
$c = $ws.Range('D47')
$c.Value = "'" + '97.40'
$c.ClearFormats()
$c = $ws.Range('E47')
$c.Value = "'" + '  -5.87%  '
$c.ClearFormats()
$c = $ws.Range('E48')
$c.Value = "'" + '  -3.25%  '
$c.ClearFormats()
$c = $ws.Range('E49')
$c.Value = "'" + '  -3.47%  '
$c.ClearFormats()
$c = $ws.Range('D50')
$c.Value = "'" + '2.18'
$c.ClearFormats()
$c = $ws.Range('E50')
$c.Value = "'" + '  -8.32%  '
$c.ClearFormats()
$c = $ws.Range('E51')
$c.Value = "'" + '  -2.72%  '
$c.ClearFormats()
